$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 134; this shifts rows 134..221 down to 135..222
$ws.Rows.Item(134).Insert()

# Populate the new row 134 with a new price-report record (week of 2022-01-21)
$ws.Range("A134").Value = 4
$ws.Range("B134").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C134").Value = "Los Lagos"
$ws.Range("D134").Value = 44582
$ws.Range("E134").Value = 10
$ws.Range("F134").Value = 100112040
$ws.Range("G134").Value = "Cilantro"
$ws.Range("H134").Value = "Sin especificar"
$ws.Range("I134").Value = "Primera"
$ws.Range("J134").Value = 120
$ws.Range("K134").Value = 12000
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = 12000
$ws.Range("N134").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O134").Value = "Región de La Araucanía"
$ws.Range("P134").Value = 6000
$ws.Range("Q134").Value = 2
$ws.Range("R134").Value = "Hortaliza"
